$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.572.96'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '2.246.21'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.48'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.642'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.24'
$ws.Range("E7").Value = '  -4.29%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.443'
$ws.Range("E9").Value = '  +2.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0960'
$ws.Range("E10").Value = '  -6.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.92'
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.66'
$ws.Range("E12").Value = '  +4.63%  '
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '2.573.29'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.46'
$ws.Range("E15").Value = '  -3.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.07'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.827'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '2.231.03'
$ws.Range("E18").Value = '  -2.19%  '
$ws.Range("D19").Value = '43.304.75'
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").Value = '0.0₃0965'
$ws.Range("E20").Value = '  -2.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.90'
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.09'
$ws.Range("E22").Value = '  -1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.51'
$ws.Range("E23").Value = '  -6.22%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.68'
$ws.Range("E25").Value = '  +31.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.41'
$ws.Range("E26").Value = '  -3.87%  '
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.74'
$ws.Range("E28").Value = '  -5.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.41'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.54'
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.131'
$ws.Range("E31").Value = '  -4.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.41'
$ws.Range("E32").Value = '  -1.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.88'
$ws.Range("E34").Value = '  +2.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0674'
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.88'
$ws.Range("E36").Value = '  -3.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.63'
$ws.Range("E37").Value = '  -6.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.29'
$ws.Range("E38").Value = '  -8.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.27'
$ws.Range("E39").Value = '  -4.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0250'
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.59'
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.52'
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.98'
$ws.Range("E44").Value = '  -3.18%  '
$ws.Range("B45").Value = 'TerraClassic'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000211'
$ws.Range("E45").Value = '  +2.38%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '96.43'
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0939'
$ws.Range("E47").Value = '  -4.27%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.17'
$ws.Range("E48").Value = '  -2.74%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.449.72'
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.05'
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.28'
$ws.Range("E51").Value = '  -3.51%  '
